# Weekly price-list update:
# A new weekly record (Fecha 45146) is inserted as row 81 of the data table,
# pushing every existing record from row 81 onward down by one row
# (old row 81 -> new row 82, ..., old row 106 -> new row 107).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 81, shifting rows 81:106 down to 82:107.
$ws.Rows("81:81").Insert()

# Populate the newly inserted row 81 with the new weekly record.
$ws.Range("A81").Value = 6
$ws.Range("B81").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C81").Value = "Metropolitana"
$ws.Range("D81").Value = 45146
$ws.Range("E81").Value = 13
$ws.Range("F81").Value = 100112035
$ws.Range("G81").Value = "Bruselas (repollito)"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 350
$ws.Range("K81").Value = 17000
$ws.Range("L81").Value = 18000
$ws.Range("M81").Value = 17343
$ws.Range("N81").Value = "$/malla 15 kilos"
$ws.Range("O81").Value = "Provincia de Quillota"
$ws.Range("P81").Value = 1156
$ws.Range("Q81").Value = 15
$ws.Range("R81").Value = "Hortaliza"

# Ensure the date cell keeps the same date number format used by the other rows.
$ws.Range("D81").NumberFormat = $ws.Range("D82").NumberFormat
